$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 325.14285
$ws.Range("I4").Value = 325.14285
$ws.Range("K4").Value = 325.14285
$ws.Range("M4").Value = -211.14285

$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H17").Value = 3057335.5
$ws.Range("J17").Value = 3057335.5
$ws.Range("L17").Value = 9172006.5
$ws.Range("N17").Value = -9172342.5

$ws.Range("H40").Value = 3891.6785
$ws.Range("I40").Value = 3736.3157
$ws.Range("J40").Value = 4219.6665
$ws.Range("K40").Value = 3736.3157
$ws.Range("L40").Value = 4219.6665
$ws.Range("M40").Value = -3561.3157
$ws.Range("N40").Value = -4569.6665

$ws.Range("H62").Value = 2249.75
$ws.Range("I62").Value = 1800.8
$ws.Range("J62").Value = 2998
$ws.Range("K62").Value = 1800.8
$ws.Range("L62").Value = 2998
$ws.Range("M62").Value = -1176.8
$ws.Range("N62").Value = -4246

$ws.Range("H65").Value = 2249.75
$ws.Range("I65").Value = 1800.8
$ws.Range("J65").Value = 2998
$ws.Range("K65").Value = 9004
$ws.Range("L65").Value = 14990
$ws.Range("M65").Value = -5884
$ws.Range("N65").Value = -21230

$ws.Range("H76").Value = 12833.429
$ws.Range("I76").Value = 12833.429
$ws.Range("K76").Value = 12833.429
$ws.Range("M76").Value = -12518.429

$ws.Range("H79").Value = 12833.429
$ws.Range("I79").Value = 12833.429
$ws.Range("K79").Value = 12833.429
$ws.Range("M79").Value = -11741.429

$ws.Range("H100").Value = 1612.1
$ws.Range("I100").Value = 1296.4375
$ws.Range("J100").Value = 2874.75
$ws.Range("K100").Value = 1296.4375
$ws.Range("L100").Value = 2874.75
$ws.Range("M100").Value = -755.4375
$ws.Range("N100").Value = -3956.75

$ws.Range("H132").Value = 10107
$ws.Range("I132").Value = 10932.7
$ws.Range("J132").Value = 1850
$ws.Range("K132").Value = 32798.10000000001
$ws.Range("L132").Value = 5550
$ws.Range("M132").Value = -30268.10000000001
$ws.Range("N132").Value = -10610

$ws.Range("H137").Value = 6632.625
$ws.Range("I137").Value = 4030.1177
$ws.Range("J137").Value = 12953
$ws.Range("K137").Value = 12090.3531
$ws.Range("L137").Value = 38859
$ws.Range("M137").Value = -9540.3531
$ws.Range("N137").Value = -43959

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 14333.667
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 21000.5
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 21000.5
$ws.Range("M30").Value = -850
$ws.Range("N30").Value = -21300.5

$ws.Range("H63").Value = 2537.75
$ws.Range("I63").Value = 2537.75
$ws.Range("K63").Value = 2537.75
$ws.Range("M63").Value = -1851.75

$ws.Range("H66").Value = 2537.75
$ws.Range("I66").Value = 2537.75
$ws.Range("K66").Value = 12688.75
$ws.Range("M66").Value = -9256.75

$ws.Range("H88").Value = 25295
$ws.Range("J88").Value = 287
$ws.Range("L88").Value = 287
$ws.Range("N88").Value = -1099

$ws.Range("H91").Value = 25295
$ws.Range("J91").Value = 287
$ws.Range("L91").Value = 287
$ws.Range("N91").Value = -3095

$ws.Range("H122").Value = 2320.96
$ws.Range("I122").Value = 2101.35
$ws.Range("K122").Value = 6304.049999999999
$ws.Range("M122").Value = -3854.049999999999

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 73770
$ws.Range("J129").Value = 73770
$ws.Range("L129").Value = 73770
$ws.Range("N129").Value = -83770

$ws.Range("H130").Value = 82141.336
$ws.Range("J130").Value = 82141.336
$ws.Range("L130").Value = 82141.336
$ws.Range("N130").Value = -92181.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1556.4348
$ws.Range("J20").Value = 1883.1666
$ws.Range("L20").Value = 1883.1666
$ws.Range("N20").Value = -2377.1666

$ws.Range("H57").Value = 47498.25
$ws.Range("I57").Value = 20000
$ws.Range("J57").Value = 74996.5
$ws.Range("K57").Value = 20000
$ws.Range("L57").Value = 74996.5
$ws.Range("M57").Value = -19280
$ws.Range("N57").Value = -76436.5

$ws.Range("H105").Value = 30964
$ws.Range("I105").Value = 30964
$ws.Range("K105").Value = 30964
$ws.Range("M105").Value = -29217

$ws.Range("H133").Value = 149990
$ws.Range("J133").Value = 149990
$ws.Range("L133").Value = 149990
$ws.Range("N133").Value = -160110

$ws.Range("H136").Value = 47498.25
$ws.Range("I136").Value = 20000
$ws.Range("J136").Value = 74996.5
$ws.Range("K136").Value = 20000
$ws.Range("L136").Value = 74996.5
$ws.Range("M136").Value = -14900
$ws.Range("N136").Value = -85196.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4352.911
$ws.Range("I31").Value = 2411.0715
$ws.Range("J31").Value = 7551.2354
$ws.Range("K31").Value = 2411.0715
$ws.Range("L31").Value = 7551.2354
$ws.Range("M31").Value = -2116.0715
$ws.Range("N31").Value = -8141.2354

$ws.Range("H34").Value = 4352.911
$ws.Range("I34").Value = 2411.0715
$ws.Range("J34").Value = 7551.2354
$ws.Range("K34").Value = 2411.0715
$ws.Range("L34").Value = 7551.2354
$ws.Range("M34").Value = -2209.0715
$ws.Range("N34").Value = -7955.2354

$ws.Range("H138").Value = 199992.25
$ws.Range("J138").Value = 199992.25
$ws.Range("L138").Value = 199992.25
$ws.Range("N138").Value = -210272.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2046.5
$ws.Range("I3").Value = 2046.5
$ws.Range("K3").Value = 6139.5
$ws.Range("M3").Value = -6027.5

$ws.Range("H5").Value = 1636.0454
$ws.Range("I5").Value = 999.6667
$ws.Range("K5").Value = 2999.0001
$ws.Range("M5").Value = -2887.0001

$ws.Range("H55").Value = 1127592.5
$ws.Range("I55").Value = 2250187.5
$ws.Range("K55").Value = 6750562.5
$ws.Range("M55").Value = -6750385.5

$ws.Range("H122").Value = 775.12195
$ws.Range("I122").Value = 415.66666
$ws.Range("J122").Value = 803.5
$ws.Range("K122").Value = 3740.99994
$ws.Range("L122").Value = 7231.5
$ws.Range("M122").Value = -1290.99994
$ws.Range("N122").Value = -12131.5

$ws.Range("H131").Value = 1451.8334
$ws.Range("J131").Value = 1541.4117
$ws.Range("L131").Value = 4624.2351
$ws.Range("N131").Value = -14704.2351

$ws.Range("H132").Value = 5314.5
$ws.Range("I132").Value = 1966
$ws.Range("J132").Value = 8663
$ws.Range("K132").Value = 17694
$ws.Range("L132").Value = 77967
$ws.Range("M132").Value = -15164
$ws.Range("N132").Value = -83027

$ws.Range("H135").Value = 1636.0454
$ws.Range("I135").Value = 999.6667
$ws.Range("K135").Value = 8997.0003
$ws.Range("M135").Value = -6462.0003

$ws.Range("H139").Value = 4117
$ws.Range("J139").Value = 3033
$ws.Range("L139").Value = 9099
$ws.Range("N139").Value = -19379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2315.6667
$ws.Range("I80").Value = 2097.5
$ws.Range("J80").Value = 2424.75
$ws.Range("K80").Value = 2097.5
$ws.Range("L80").Value = 2424.75
$ws.Range("M80").Value = -1099.5
$ws.Range("N80").Value = -4420.75

$ws.Range("H83").Value = 2315.6667
$ws.Range("I83").Value = 2097.5
$ws.Range("J83").Value = 2424.75
$ws.Range("K83").Value = 10487.5
$ws.Range("L83").Value = 12123.75
$ws.Range("M83").Value = -5495.5
$ws.Range("N83").Value = -22107.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2655.7
$ws.Range("I16").Value = 2255.6956
$ws.Range("K16").Value = 2255.6956
$ws.Range("M16").Value = -2085.6956

$ws.Range("H17").Value = 7500
$ws.Range("I17").Value = 7500
$ws.Range("K17").Value = 7500
$ws.Range("M17").Value = -7330

$ws.Range("H22").Value = 2162.6667
$ws.Range("I22").Value = 2101.4614
$ws.Range("J22").Value = 2321.8
$ws.Range("K22").Value = 2101.4614
$ws.Range("L22").Value = 2321.8
$ws.Range("M22").Value = -1806.4614
$ws.Range("N22").Value = -2911.8

$ws.Range("H27").Value = 2162.6667
$ws.Range("I27").Value = 2101.4614
$ws.Range("J27").Value = 2321.8
$ws.Range("K27").Value = 2101.4614
$ws.Range("L27").Value = 2321.8
$ws.Range("M27").Value = -1994.4614
$ws.Range("N27").Value = -2535.8

$ws.Range("H122").Value = 3279.3635
$ws.Range("I122").Value = 2830.3333
$ws.Range("J122").Value = 5300
$ws.Range("K122").Value = 8490.999899999999
$ws.Range("L122").Value = 15900
$ws.Range("M122").Value = -6040.999899999999
$ws.Range("N122").Value = -20800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 17999.5
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 17999.5
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 17999.5
$ws.Range("N24").Value = -18459.5
$ws.Range("M24").ClearContents()
